$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '72.391.69'
$ws.Range("E2").Value = '  +0.18%  '
$ws.Range("D3").Value = '2.635.93'
$ws.Range("E3").Value = '  -1.33%  '
$ws.Range("E4").Value = '  +0.02%  '
$ws.Range("D5").Value = '585.18'
$ws.Range("E5").Value = '  -2.36%  '
$ws.Range("D6").Value = '175.21'
$ws.Range("E6").Value = '  -0.58%  '
$ws.Range("E7").Value = '  +0.01%  '
$ws.Range("D8").Value = '0.519'
$ws.Range("E8").Value = '  -0.63%  '
$ws.Range("D9").Value = '2.636.39'
$ws.Range("E9").Value = '  -1.27%  '
$ws.Range("D10").Value = '0.171'
$ws.Range("E10").Value = '  +1.43%  '
$ws.Range("E11").Value = '  +1.43%  '
$ws.Range("D12").Value = '0.359'
$ws.Range("E12").Value = '  +1.64%  '
$ws.Range("D13").Value = '4.93'
$ws.Range("E13").Value = '  -1.82%  '
$ws.Range("D14").Value = '3.121.06'
$ws.Range("E14").Value = '  -1.38%  '
$ws.Range("E15").Value = '  +0.42%  '
$ws.Range("D16").Value = '72.196.74'
$ws.Range("E16").Value = '  +0.13%  '
$ws.Range("D17").Value = '25.72'
$ws.Range("E17").Value = '  -2.07%  '
$ws.Range("D18").Value = '2.666.06'
$ws.Range("E18").Value = '  -0.56%  '
$ws.Range("D19").Value = '12.03'
$ws.Range("E19").Value = '  +0.02%  '
$ws.Range("D20").Value = '376.17'
$ws.Range("E20").Value = '  +1.66%  '
$ws.Range("D21").Value = '7.85'
$ws.Range("E21").Value = '  -1.63%  '
$origStyle = $ws.Range("D22").Style
$ws.Range("D22").Value = "'4.10"
$ws.Range("D22").Style = $origStyle
$ws.Range("E22").Value = '  -1.22%  '
$ws.Range("E23").Value = '  -0.63%  '
$origStyle = $ws.Range("D24").Style
$ws.Range("D24").Value = "'71.50"
$ws.Range("D24").Style = $origStyle
$ws.Range("E24").Value = '  -0.37%  '
$ws.Range("E25").Value = '  +0.00%  '
$ws.Range("D26").Value = '4.23'
$ws.Range("E26").Value = '  -2.43%  '
$ws.Range("D27").Value = '9.51'
$ws.Range("E27").Value = '  -3.02%  '
$ws.Range("D28").Value = '2.775.30'
$ws.Range("E28").Value = '  -1.19%  '
$ws.Range("E29").Value = '  -0.26%  '
$ws.Range("D30").Value = '0.0₃0948'
$ws.Range("E30").Value = '  +1.00%  '
$ws.Range("D31").Value = '7.96'
$ws.Range("E31").Value = '  -1.17%  '
$ws.Range("D32").Value = '491.27'
$ws.Range("D33").Value = '1.31'
$ws.Range("E33").Value = '  +1.78%  '
$ws.Range("D34").Value = '1.79'
$ws.Range("E34").Value = '  -0.96%  '
$ws.Range("E35").Value = '  +0.02%  '
$ws.Range("D36").Value = '161.67'
$ws.Range("E36").Value = '  -1.89%  '
$ws.Range("D37").Value = '0.115'
$ws.Range("E37").Value = '  +9.17%  '
$ws.Range("D38").Value = '19.16'
$ws.Range("E38").Value = '  -2.01%  '
$ws.Range("D39").Value = '18.89'
$ws.Range("E39").Value = '  -1.08%  '
$ws.Range("E40").Value = '  -1.10%  '
$ws.Range("E41").Value = '  +0.00%  '
$ws.Range("D42").Value = '1.72'
$ws.Range("E42").Value = '  -5.12%  '
$ws.Range("D43").Value = '2.57'
$ws.Range("E43").Value = '  -0.03%  '
$ws.Range("D44").Value = '4.88'
$ws.Range("E44").Value = '  -2.45%  '
$ws.Range("D45").Value = '0.325'
$ws.Range("E45").Value = '  -2.05%  '
$ws.Range("D46").Value = '39.03'
$ws.Range("E46").Value = '  -0.57%  '
$ws.Range("D47").Value = '150.24'
$ws.Range("E47").Value = '  -2.07%  '
$ws.Range("B48").Value = 'Filecoin'
$ws.Range("C48").Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range("D48").Value = '3.62'
$ws.Range("E48").Value = '  -2.60%  '
$ws.Range("B49").Value = 'ARBITRUM'
$ws.Range("C49").Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range("D49").Value = '0.541'
$ws.Range("E49").Value = '  -1.27%  '
$ws.Range("D50").Value = '1.67'
$ws.Range("E50").Value = '  -2.74%  '
$ws.Range("D51").Value = '0.608'
$ws.Range("E51").Value = '  +1.21%  '
